$d = $word.ActiveDocument
$footer = $d.Sections(1).Footers(1)
$r = $footer.Range.Duplicate
$r.Collapse(0)
$r.InsertAfter(".1")
$r.Font.Name = "Cambria"
$r.Font.Size = 8
$r.Font.Italic = $true
$r.NoProofing = $true
Write-Output $footer.Range.Text
